$d = $word.ActiveDocument

function Merge-Text($oldText) {
    # Re-applying the same text over a Find match causes the interop
    # engine to rewrite the matched span as a single run, which is how
    # the previously split runs (identical formatting) get coalesced.
    $null = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)
}

Merge-Text "Declaraciones: define todas las posibles declaraciones que pueden darse en el lenguaje."
Merge-Text "Tipo: esta clase representa todos los tipos posibles de las variables definidas en el lenguaje (int, real o bool)."
Merge-Text "Carácter: esta clase representa todos los caracteres posibles que pueden aparecer en nuestro lenguaje (a-z, 0-9, _)."
Merge-Text "Asignación: esta clase genera todas las posibles asignaciones que pueden existir, en nuestro caso solo asignaciones variable = expresión."
Merge-Text "DígitoPositivo: genera un dígito positivo (1-9)."
Merge-Text "ExpCompleja: genera todo el conjunto posible de expresiones del lenguaje."
Merge-Text "DUDA:"
Merge-Text "Definiciones auxiliares"
Merge-Text "Programa ->Sec_Declaraciones && Sec_Instrucciones"
Merge-Text "Variable -> Letra(Carácter*) "
Merge-Text "Decimal -> .(Dígito*)DígitoPositivo"

# --- NOTA paragraph: insert "(blancos)" after "NF" and relocate the
#     _GoBack bookmark to sit right after the newly inserted text. ---

$r = $d.Content
$null = $r.Find.Execute("NF", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(" (blancos)")

# Re-find the text we just inserted so we can bracket it with a
# throwaway bookmark -- inserting a bookmark around a run forces the
# engine to keep it as its own distinct run instead of folding it back
# into its identically-formatted neighbour.
$r2 = $d.Content
$null = $r2.Find.Execute(" (blancos)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("TEMP_SPLIT", $r2)

# Collapse to the point right after "(blancos)" and drop the _GoBack
# bookmark there (this also removes it from its old location next to
# the diagram, since a document can only have one bookmark per name).
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2)

# Clean up the temporary helper bookmark.
$bm = $d.Bookmarks("TEMP_SPLIT")
$bm.Delete()
